$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in A1
$ws.Cells.Item(1,1).Value = "Datos actualizados a 6 de Mayo de 2020 a las 20:03"

# Row 4: Estados Unidos
$ws.Cells.Item(4,1).Value = "Estados Unidos"
$ws.Cells.Item(4,2).Value = 1245857
$ws.Cells.Item(4,3).Value = 8224
$ws.Cells.Item(4,4).Value = 203088
$ws.Cells.Item(4,5).Value = 969624
$ws.Cells.Item(4,6).Value = 16173
$ws.Cells.Item(4,7).Value = 874
$ws.Cells.Item(4,8).Value = 73145

# Row 9: Alemania
$ws.Cells.Item(9,1).Value = "Alemania"
$ws.Cells.Item(9,2).Value = 167372
$ws.Cells.Item(9,3).Value = 365
$ws.Cells.Item(9,4).Value = 137400
$ws.Cells.Item(9,5).Value = 22979
$ws.Cells.Item(9,6).Value = 1884
$ws.Cells.Item(9,7).Value = 0
$ws.Cells.Item(9,8).Value = 6993

# Row 15: Canada
$ws.Cells.Item(15,1).Value = "Canada"
$ws.Cells.Item(15,2).Value = 63375
$ws.Cells.Item(15,3).Value = 1329
$ws.Cells.Item(15,4).Value = 26993
$ws.Cells.Item(15,5).Value = 32159
$ws.Cells.Item(15,6).Value = 502
$ws.Cells.Item(15,7).Value = 180
$ws.Cells.Item(15,8).Value = 4223

# Row 33: Emiratos Arabes Unidos
$ws.Cells.Item(33,1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(33,2).Value = 15738
$ws.Cells.Item(33,3).Value = 546
$ws.Cells.Item(33,4).Value = 3359
$ws.Cells.Item(33,5).Value = 12222
$ws.Cells.Item(33,6).Value = 1
$ws.Cells.Item(33,7).Value = 11
$ws.Cells.Item(33,8).Value = 157

# Row 34: Austria
$ws.Cells.Item(34,1).Value = "Austria"
$ws.Cells.Item(34,2).Value = 15684
$ws.Cells.Item(34,3).Value = 34
$ws.Cells.Item(34,4).Value = 13639
$ws.Cells.Item(34,5).Value = 1437
$ws.Cells.Item(34,6).Value = 97
$ws.Cells.Item(34,7).Value = 2
$ws.Cells.Item(34,8).Value = 608

# Row 35: Japon
$ws.Cells.Item(35,1).Value = "Japon"
$ws.Cells.Item(35,2).Value = 15253
$ws.Cells.Item(35,3).Value = 0
$ws.Cells.Item(35,4).Value = 4496
$ws.Cells.Item(35,5).Value = 10201
$ws.Cells.Item(35,6).Value = 308
$ws.Cells.Item(35,7).Value = 0
$ws.Cells.Item(35,8).Value = 556

# Row 37: Rumania
$ws.Cells.Item(37,1).Value = "Rumania"
$ws.Cells.Item(37,2).Value = 14107
$ws.Cells.Item(37,3).Value = 270
$ws.Cells.Item(37,4).Value = 5788
$ws.Cells.Item(37,5).Value = 7455
$ws.Cells.Item(37,6).Value = 244
$ws.Cells.Item(37,7).Value = 23
$ws.Cells.Item(37,8).Value = 864

# Row 60: Kazajistan
$ws.Cells.Item(60,1).Value = "Kazajistan"
$ws.Cells.Item(60,2).Value = 4344
$ws.Cells.Item(60,3).Value = 139
$ws.Cells.Item(60,4).Value = 1408
$ws.Cells.Item(60,5).Value = 2906
$ws.Cells.Item(60,6).Value = 31
$ws.Cells.Item(60,7).Value = 1
$ws.Cells.Item(60,8).Value = 30

# Row 61: Barein
$ws.Cells.Item(61,1).Value = "Barein"
$ws.Cells.Item(61,2).Value = 3934
$ws.Cells.Item(61,3).Value = 214
$ws.Cells.Item(61,4).Value = 1860
$ws.Cells.Item(61,5).Value = 2066
$ws.Cells.Item(61,6).Value = 4
$ws.Cells.Item(61,7).Value = 0
$ws.Cells.Item(61,8).Value = 8

# Row 62: Luxemburgo
$ws.Cells.Item(62,1).Value = "Luxemburgo"
$ws.Cells.Item(62,2).Value = 3851
$ws.Cells.Item(62,3).Value = 11
$ws.Cells.Item(62,4).Value = 3452
$ws.Cells.Item(62,5).Value = 301
$ws.Cells.Item(62,6).Value = 20
$ws.Cells.Item(62,7).Value = 2
$ws.Cells.Item(62,8).Value = 98

# Row 77: Guinea
$ws.Cells.Item(77,1).Value = "Guinea"
$ws.Cells.Item(77,2).Value = 1856
$ws.Cells.Item(77,3).Value = 45
$ws.Cells.Item(77,4).Value = 597
$ws.Cells.Item(77,5).Value = 1248
$ws.Cells.Item(77,6).Value = 0
$ws.Cells.Item(77,7).Value = 1
$ws.Cells.Item(77,8).Value = 11

# Row 105: Principado de Andorra
$ws.Cells.Item(105,1).Value = "Principado de Andorra"
$ws.Cells.Item(105,2).Value = 751
$ws.Cells.Item(105,3).Value = 0
$ws.Cells.Item(105,4).Value = 521
$ws.Cells.Item(105,5).Value = 184
$ws.Cells.Item(105,6).Value = 15
$ws.Cells.Item(105,7).Value = 0
$ws.Cells.Item(105,8).Value = 46

# Row 108: Burkina Faso
$ws.Cells.Item(108,1).Value = "Burkina Faso"
$ws.Cells.Item(108,2).Value = 729
$ws.Cells.Item(108,3).Value = 41
$ws.Cells.Item(108,4).Value = 555
$ws.Cells.Item(108,5).Value = 126
$ws.Cells.Item(108,6).Value = 0
$ws.Cells.Item(108,7).Value = 0
$ws.Cells.Item(108,8).Value = 48

# Row 109: Crucero
$ws.Cells.Item(109,1).Value = "Crucero"
$ws.Cells.Item(109,2).Value = 712
$ws.Cells.Item(109,3).Value = 0
$ws.Cells.Item(109,4).Value = 645
$ws.Cells.Item(109,5).Value = 54
$ws.Cells.Item(109,6).Value = 4
$ws.Cells.Item(109,7).Value = 0
$ws.Cells.Item(109,8).Value = 13

# Row 111: El Salvador
$ws.Cells.Item(111,1).Value = "El Salvador"
$ws.Cells.Item(111,2).Value = 633
$ws.Cells.Item(111,3).Value = 46
$ws.Cells.Item(111,4).Value = 219
$ws.Cells.Item(111,5).Value = 399
$ws.Cells.Item(111,6).Value = 4
$ws.Cells.Item(111,7).Value = 2
$ws.Cells.Item(111,8).Value = 15

# Row 122: Paraguay
$ws.Cells.Item(122,1).Value = "Paraguay"
$ws.Cells.Item(122,2).Value = 440
$ws.Cells.Item(122,3).Value = 9
$ws.Cells.Item(122,4).Value = 142
$ws.Cells.Item(122,5).Value = 288
$ws.Cells.Item(122,6).Value = 9
$ws.Cells.Item(122,7).Value = 0
$ws.Cells.Item(122,8).Value = 10

# Row 123: Guinea Ecuatorial
$ws.Cells.Item(123,1).Value = "Guinea Ecuatorial"
$ws.Cells.Item(123,2).Value = 439
$ws.Cells.Item(123,3).Value = 124
$ws.Cells.Item(123,4).Value = 13
$ws.Cells.Item(123,5).Value = 422
$ws.Cells.Item(123,6).Value = 0
$ws.Cells.Item(123,7).Value = 1
$ws.Cells.Item(123,8).Value = 4

# Row 124: Taiwan
$ws.Cells.Item(124,1).Value = "Taiwan"
$ws.Cells.Item(124,2).Value = 439
$ws.Cells.Item(124,3).Value = 1
$ws.Cells.Item(124,4).Value = 339
$ws.Cells.Item(124,5).Value = 94
$ws.Cells.Item(124,6).Value = 0
$ws.Cells.Item(124,7).Value = 0
$ws.Cells.Item(124,8).Value = 6

# Row 152: Suazilandia
$ws.Cells.Item(152,1).Value = "Suazilandia"
$ws.Cells.Item(152,2).Value = 123
$ws.Cells.Item(152,3).Value = 4
$ws.Cells.Item(152,4).Value = 12
$ws.Cells.Item(152,5).Value = 109
$ws.Cells.Item(152,6).Value = 0
$ws.Cells.Item(152,7).Value = 1
$ws.Cells.Item(152,8).Value = 2

# Row 153: Camboya
$ws.Cells.Item(153,1).Value = "Camboya"
$ws.Cells.Item(153,2).Value = 122
$ws.Cells.Item(153,3).Value = 0
$ws.Cells.Item(153,4).Value = 120
$ws.Cells.Item(153,5).Value = 2
$ws.Cells.Item(153,6).Value = 1
$ws.Cells.Item(153,7).Value = 0
$ws.Cells.Item(153,8).Value = 0

# Row 157: Aruba
$ws.Cells.Item(157,1).Value = "Aruba"
$ws.Cells.Item(157,2).Value = 101
$ws.Cells.Item(157,3).Value = 0
$ws.Cells.Item(157,4).Value = 89
$ws.Cells.Item(157,5).Value = 10
$ws.Cells.Item(157,6).Value = 4
$ws.Cells.Item(157,7).Value = 0
$ws.Cells.Item(157,8).Value = 2

# Row 205: Seychelles
$ws.Cells.Item(205,1).Value = "Seychelles"
$ws.Cells.Item(205,2).Value = 11
$ws.Cells.Item(205,3).Value = 0
$ws.Cells.Item(205,4).Value = 8
$ws.Cells.Item(205,5).Value = 3
$ws.Cells.Item(205,6).Value = 0
$ws.Cells.Item(205,7).Value = 0
$ws.Cells.Item(205,8).Value = 0

# Row 206: Montserrat
$ws.Cells.Item(206,1).Value = "Montserrat"
$ws.Cells.Item(206,2).Value = 11
$ws.Cells.Item(206,3).Value = 0
$ws.Cells.Item(206,4).Value = 7
$ws.Cells.Item(206,5).Value = 3
$ws.Cells.Item(206,6).Value = 1
$ws.Cells.Item(206,7).Value = 0
$ws.Cells.Item(206,8).Value = 1
